$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 278.57144
$ws.Range("I9").Value = 250
$ws.Range("K9").Value = 250
$ws.Range("M9").Value = -81
$ws.Range("H12").Value = 1010.5
$ws.Range("I12").Value = 228
$ws.Range("K12").Value = 228
$ws.Range("M12").Value = -58
$ws.Range("H19").Value = 1835.3
$ws.Range("I19").Value = 799.3333
$ws.Range("K19").Value = 799.3333
$ws.Range("M19").Value = -624.3333
$ws.Range("H28").Value = 2039.1538
$ws.Range("I28").Value = 2266.5
$ws.Range("J28").Value = 1675.4
$ws.Range("K28").Value = 2266.5
$ws.Range("L28").Value = 1675.4
$ws.Range("M28").Value = -1781.5
$ws.Range("N28").Value = -2645.4
$ws.Range("H86").Value = 3536
$ws.Range("I86").Value = 3249.5
$ws.Range("J86").Value = 4300
$ws.Range("K86").Value = 3249.5
$ws.Range("L86").Value = 4300
$ws.Range("M86").Value = -2126.5
$ws.Range("N86").Value = -6546
$ws.Range("H89").Value = 3536
$ws.Range("I89").Value = 3249.5
$ws.Range("J89").Value = 4300
$ws.Range("K89").Value = 16247.5
$ws.Range("L89").Value = 21500
$ws.Range("M89").Value = -10631.5
$ws.Range("N89").Value = -32732
$ws.Range("H116").Value = 4452.8
$ws.Range("I116").Value = 3819.2
$ws.Range("K116").Value = 3819.2
$ws.Range("M116").Value = -377.1999999999998
$ws.Range("H137").Value = 1712.0264
$ws.Range("I137").Value = 1625.2667
$ws.Range("J137").Value = 2037.375
$ws.Range("K137").Value = 4875.800099999999
$ws.Range("L137").Value = 6112.125
$ws.Range("M137").Value = -2325.800099999999
$ws.Range("N137").Value = -11212.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1748.8462
$ws.Range("I2").Value = 1061.25
$ws.Range("K2").Value = 1061.25
$ws.Range("M2").Value = -948.25
$ws.Range("H32").Value = 5203.518
$ws.Range("I32").Value = 5207.2363
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 5207.2363
$ws.Range("L32").Value = 4999
$ws.Range("M32").Value = -4920.2363
$ws.Range("N32").Value = -5573
$ws.Range("H116").Value = 1748.8462
$ws.Range("I116").Value = 1061.25
$ws.Range("K116").Value = 1061.25
$ws.Range("M116").Value = 1232.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1748.8462
$ws.Range("I3").Value = 1061.25
$ws.Range("K3").Value = 1061.25
$ws.Range("M3").Value = -947.25
$ws.Range("H22").Value = 2024.1364
$ws.Range("J22").Value = 950
$ws.Range("L22").Value = 950
$ws.Range("N22").Value = -1296
$ws.Range("H54").Value = 4499.75
$ws.Range("I54").Value = 4499.75
$ws.Range("K54").Value = 4499.75
$ws.Range("M54").Value = -4015.75
$ws.Range("H82").Value = 24475
$ws.Range("I82").Value = 3950
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 3950
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -3567
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 24475
$ws.Range("I85").Value = 3950
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 3950
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -2624
$ws.Range("N85").Value = -47652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3798.2068
$ws.Range("I31").Value = 3969.3
$ws.Range("J31").Value = 3708.158
$ws.Range("K31").Value = 3969.3
$ws.Range("L31").Value = 3708.158
$ws.Range("M31").Value = -3674.3
$ws.Range("N31").Value = -4298.157999999999
$ws.Range("H34").Value = 3798.2068
$ws.Range("I34").Value = 3969.3
$ws.Range("J34").Value = 3708.158
$ws.Range("K34").Value = 3969.3
$ws.Range("L34").Value = 3708.158
$ws.Range("M34").Value = -3767.3
$ws.Range("N34").Value = -4112.157999999999
$ws.Range("H99").Value = 5050.6665
$ws.Range("I99").Value = 5261
$ws.Range("J99").Value = 3999
$ws.Range("K99").Value = 5261
$ws.Range("L99").Value = 3999
$ws.Range("M99").Value = -3763
$ws.Range("N99").Value = -6995
$ws.Range("H122").Value = 3950.875
$ws.Range("I122").Value = 3931.2
$ws.Range("K122").Value = 11793.6
$ws.Range("M122").Value = -9343.599999999999
$ws.Range("H126").Value = 5050.6665
$ws.Range("I126").Value = 5261
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 15783
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -13313
$ws.Range("N126").Value = -16937

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20487538
$ws.Range("I4").Value = 15543372
$ws.Range("J4").Value = 30375870
$ws.Range("K4").Value = 46630116
$ws.Range("L4").Value = 91127610
$ws.Range("M4").Value = -46630004
$ws.Range("N4").Value = -91127834
$ws.Range("H11").Value = 42000000
$ws.Range("I11").Value = 42000000
$ws.Range("K11").Value = 126000000
$ws.Range("M11").Value = -125999860
$ws.Range("H23").Value = 193
$ws.Range("I23").Value = 245
$ws.Range("J23").Value = 167
$ws.Range("K23").Value = 735
$ws.Range("L23").Value = 501
$ws.Range("M23").Value = -500
$ws.Range("N23").Value = -971
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("K56").Value = 15000
$ws.Range("M56").Value = -14470
$ws.Range("H112").Value = 4198.8
$ws.Range("H121").Value = 10101055
$ws.Range("I121").Value = 704.8
$ws.Range("J121").Value = 20201406
$ws.Range("K121").Value = 2114.4
$ws.Range("L121").Value = 60604218
$ws.Range("M121").Value = -804.3999999999996
$ws.Range("N121").Value = -60606838
$ws.Range("H122").Value = 501.94116
$ws.Range("J122").Value = 696.3
$ws.Range("L122").Value = 6266.7
$ws.Range("N122").Value = -11166.7
$ws.Range("H124").Value = 1370
$ws.Range("J124").Value = 2166.6667
$ws.Range("L124").Value = 6500.000100000001
$ws.Range("N124").Value = -16320.0001
$ws.Range("H133").Value = 2984.5
$ws.Range("I133").Value = 2984.5
$ws.Range("K133").Value = 8953.5
$ws.Range("M133").Value = -3893.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16732.666
$ws.Range("J70").Value = 17555.111
$ws.Range("L70").Value = 17555.111
$ws.Range("N70").Value = -18095.111
$ws.Range("H73").Value = 16732.666
$ws.Range("J73").Value = 17555.111
$ws.Range("L73").Value = 17555.111
$ws.Range("N73").Value = -19427.111
$ws.Range("H97").Value = 884.34784
$ws.Range("I97").Value = 778.0952
$ws.Range("K97").Value = 778.0952
$ws.Range("M97").Value = -282.0952
$ws.Range("H102").Value = 5103.4165
$ws.Range("I102").Value = 3659
$ws.Range("J102").Value = 6135.143
$ws.Range("K102").Value = 3659
$ws.Range("L102").Value = 6135.143
$ws.Range("M102").Value = -2037
$ws.Range("N102").Value = -9379.143
$ws.Range("H113").Value = 225506.56
$ws.Range("I113").Value = 366376.72
$ws.Range("J113").Value = 4139.143
$ws.Range("K113").Value = 366376.72
$ws.Range("L113").Value = 4139.143
$ws.Range("M113").Value = -364206.72
$ws.Range("N113").Value = -8479.143
$ws.Range("H126").Value = 5465.8887
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1645.7858
$ws.Range("I16").Value = 1745.4
$ws.Range("J16").Value = 1396.75
$ws.Range("K16").Value = 1745.4
$ws.Range("L16").Value = 1396.75
$ws.Range("M16").Value = -1575.4
$ws.Range("N16").Value = -1736.75
$ws.Range("H55").Value = 501.45456
$ws.Range("J55").Value = 414.14285
$ws.Range("L55").Value = 414.14285
$ws.Range("N55").Value = -760.14285
$ws.Range("H82").Value = 1313.8889
$ws.Range("J82").Value = 1656
$ws.Range("L82").Value = 1656
$ws.Range("N82").Value = -2378
$ws.Range("H85").Value = 1313.8889
$ws.Range("J85").Value = 1656
$ws.Range("L85").Value = 1656
$ws.Range("N85").Value = -4152
$ws.Range("H122").Value = 4364.0625
$ws.Range("I122").Value = 3728.625
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 11185.875
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -8735.875
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 8997.130999999999
$ws.Range("I132").Value = 9021.755999999999
$ws.Range("K132").Value = 27065.268
$ws.Range("M132").Value = -24535.268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9996.25
$ws.Range("J15").Value = 9995
$ws.Range("L15").Value = 9995
$ws.Range("N15").Value = -10571
$ws.Range("H74").Value = 21200.428
$ws.Range("J74").Value = 21200.428
$ws.Range("L74").Value = 21200.428
$ws.Range("N74").Value = -23072.428
$ws.Range("H77").Value = 21200.428
$ws.Range("J77").Value = 21200.428
$ws.Range("L77").Value = 63601.284
$ws.Range("N77").Value = -72961.284
$ws.Range("H122").Value = 5997.5
$ws.Range("I122").Value = 5568.357
$ws.Range("J122").Value = 6998.8335
$ws.Range("K122").Value = 16705.071
$ws.Range("L122").Value = 20996.5005
$ws.Range("M122").Value = -14255.071
$ws.Range("N122").Value = -25896.5005
